# Update the "Data" worksheet with the refreshed MPA test automation
# upload identifiers / sequence numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("K6").Value  = 60000427
$ws.Range("L7").Value  = 292
$ws.Range("N7").Value  = 60000428
$ws.Range("L8").Value  = 292
$ws.Range("O8").Value  = 293
$ws.Range("K9").Value  = 60000427
$ws.Range("L10").Value = 292
$ws.Range("K11").Value = 60000427
$ws.Range("N11").Value = 60000428
$ws.Range("L12").Value = 292
$ws.Range("N12").Value = 60000428
$ws.Range("L13").Value = 292
$ws.Range("O13").Value = 293
$ws.Range("K14").Value = 60000427
$ws.Range("L15").Value = 292
$ws.Range("K16").Value = 60000427
$ws.Range("N16").Value = 60000428
$ws.Range("L17").Value = 292
$ws.Range("N17").Value = 60000428
$ws.Range("L18").Value = 292
$ws.Range("O18").Value = 293
$ws.Range("K19").Value = 60000427
$ws.Range("L20").Value = 292
$ws.Range("K21").Value = 60000427
$ws.Range("N21").Value = 60000428
$ws.Range("L22").Value = 292
$ws.Range("N22").Value = 60000428
$ws.Range("L23").Value = 292
$ws.Range("O23").Value = 293
$ws.Range("K24").Value = 60000427
$ws.Range("L25").Value = 292
$ws.Range("K26").Value = 60000427
$ws.Range("N26").Value = 60000428
$ws.Range("L27").Value = 292
$ws.Range("N27").Value = 60000428
$ws.Range("L28").Value = 292
$ws.Range("O28").Value = 293
$ws.Range("K29").Value = 60000427
